# daily auto push: 2026-02-13 19:06 UTC
# Insert a new daily record (2026/02/14, 土, 1, 201) before the existing
# 2026/12/29 block, shifting rows 816:857 down to 817:858.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row at 816; everything currently at/after row 816
# (the 2026/12/29 ... 2027/01/05 data) shifts down by one row.
$ws.Rows("816:816").Insert()

# Column A holds the date formatted as plain text (e.g. "2026/02/13"),
# not a real Excel date. Temporarily force a text number format before
# assigning the value so Excel doesn't auto-convert the string into a
# date serial number, then clear the format again so the new cell keeps
# the same (default/general) style as every other data row.
$ws.Range("A816").NumberFormat = "@"
$ws.Range("A816").Value = "2026/02/14"
$ws.Range("A816").ClearFormats()

$ws.Range("B816").Value = "土"
$ws.Range("C816").Value = 1
$ws.Range("D816").Value = 201
